# HRM-mirror simulation parameter optimization
# Update the grating pitch / line density input on the "HRM" sheet from 20 to 100.
# All other cells on this sheet are formulas that depend (directly or
# transitively) on B3, so Excel recalculates them automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HRM")

$ws.Range("B3").Value = 100

# Leave the cursor/selection where the author left it when saving.
$ws.Range("L11").Select() | Out-Null
